# "Add evaluator and scripts"
# The evaluator script now also emits a confusion-matrix-sums table on the
# "Projects" sheet, and the original SAD-SAM-Code gold-standard-link column
# is replaced by a SAD-Code column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Projects")

# --- New "Confusion matrix sums" block (rows 11-17) -----------------------
$ws.Range("B11").Value = "Confusion matrix sums"

$ws.Range("B12").Value = "Project"
$ws.Range("C12").Value = "SAD-Code"
$ws.Range("D12").Value = "SAD-SAM"
$ws.Range("E12").Value = "SAM-Code"

$ws.Range("B13").Value = "MediaStore"
$ws.Range("C13").Value = 3589
$ws.Range("D13").Value = 518
$ws.Range("E13").Value = 2231

$ws.Range("B14").Value = "TeaStore"
$ws.Range("C14").Value = 8815
$ws.Range("D14").Value = 473
$ws.Range("E14").Value = 3895

$ws.Range("B15").Value = "TEAMMATES"
$ws.Range("C15").Value = 165330
$ws.Range("D15").Value = 1584
$ws.Range("E15").Value = 13360

$ws.Range("B16").Value = "BigBlueButton"
$ws.Range("C16").Value = 47600
$ws.Range("D16").Value = 1020
$ws.Range("E16").Value = 13440

$ws.Range("B17").Value = "JabRef"
$ws.Range("C17").Value = 26000
$ws.Range("D17").Value = 78
$ws.Range("E17").Value = 12000

# --- Rename the original "SAD-SAM-Code" header to "SAD-Code" --------------
$ws.Range("C3").Value = "SAD-Code"

# --- Make "Projects" the active sheet/selection, matching the saved file --
$ws.Range("D17").Select()
$ws.Activate()
